# Insert a new data row above the current row 290 (Excel row numbers), which
# shifts all existing rows 290-386 down to 291-387, then populate the newly
# inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 290 - pushes existing row 290 (and everything below)
# down by one row.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new record's data.
$ws.Cells.Item(290, 1).Value = 5
$ws.Cells.Item(290, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(290, 3).Value = "Maule"
$ws.Cells.Item(290, 4).Value = 44524
$ws.Cells.Item(290, 5).Value = 7
$ws.Cells.Item(290, 6).Value = 100112002
$ws.Cells.Item(290, 7).Value = "Pimiento"
$ws.Cells.Item(290, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 200
$ws.Cells.Item(290, 11).Value = 15000
$ws.Cells.Item(290, 12).Value = 15000
$ws.Cells.Item(290, 13).Value = 15000
$ws.Cells.Item(290, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(290, 15).Value = "Región del Maule"
$ws.Cells.Item(290, 16).Value = 1000
$ws.Cells.Item(290, 17).Value = 15
$ws.Cells.Item(290, 18).Value = "Hortaliza"
